$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9536.773999999999
$ws.Range("I62").Value = 10044.1875
$ws.Range("K62").Value = 10044.1875
$ws.Range("M62").Value = -9420.1875
$ws.Range("H65").Value = 9536.773999999999
$ws.Range("I65").Value = 10044.1875
$ws.Range("K65").Value = 50220.9375
$ws.Range("M65").Value = -47100.9375
$ws.Range("H69").Value = 14887.556
$ws.Range("I69").Value = 10498.25
$ws.Range("J69").Value = 18399
$ws.Range("K69").Value = 31494.75
$ws.Range("L69").Value = 55197
$ws.Range("M69").Value = -30620.75
$ws.Range("N69").Value = -56945
$ws.Range("H72").Value = 14887.556
$ws.Range("I72").Value = 10498.25
$ws.Range("J72").Value = 18399
$ws.Range("K72").Value = 94484.25
$ws.Range("L72").Value = 165591
$ws.Range("M72").Value = -90116.25
$ws.Range("N72").Value = -174327
$ws.Range("H92").Value = 144249.58
$ws.Range("J92").Value = 1629
$ws.Range("L92").Value = 1629
$ws.Range("N92").Value = -4125
$ws.Range("H103").Value = 1654.1666
$ws.Range("I103").Value = 1784.7142
$ws.Range("J103").Value = 1471.4
$ws.Range("K103").Value = 5354.142599999999
$ws.Range("L103").Value = 4414.200000000001
$ws.Range("M103").Value = -4768.142599999999
$ws.Range("N103").Value = -5586.200000000001
$ws.Range("H138").Value = 6196.4
$ws.Range("I138").Value = 5397.1
$ws.Range("J138").Value = 6995.7
$ws.Range("K138").Value = 16191.3
$ws.Range("L138").Value = 20987.1
$ws.Range("M138").Value = -11051.3
$ws.Range("N138").Value = -31267.1
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 111111
$ws.Range("I26").Value = 111111
$ws.Range("K26").Value = 111111
$ws.Range("M26").Value = -110781
$ws.Range("H32").Value = 27948.143
$ws.Range("I32").Value = 31011.771
$ws.Range("K32").Value = 31011.771
$ws.Range("M32").Value = -30724.771
$ws.Range("H45").Value = 2780.818
$ws.Range("I45").Value = 979.1818
$ws.Range("J45").Value = 4582.4546
$ws.Range("K45").Value = 979.1818
$ws.Range("L45").Value = 4582.4546
$ws.Range("M45").Value = -602.1818
$ws.Range("N45").Value = -5336.4546
$ws.Range("H107").Value = 149499
$ws.Range("J107").Value = 149499
$ws.Range("L107").Value = 149499
$ws.Range("N107").Value = -157179
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 111288.11
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H80").Value = 474.81818
$ws.Range("J80").Value = 356.25
$ws.Range("L80").Value = 356.25
$ws.Range("N80").Value = -2352.25
$ws.Range("H81").Value = 209023.17
$ws.Range("I81").Value = 13709
$ws.Range("J81").Value = 248086
$ws.Range("K81").Value = 13709
$ws.Range("L81").Value = 248086
$ws.Range("M81").Value = -12648
$ws.Range("N81").Value = -250208
$ws.Range("H83").Value = 474.81818
$ws.Range("J83").Value = 356.25
$ws.Range("L83").Value = 1781.25
$ws.Range("N83").Value = -11765.25
$ws.Range("H84").Value = 209023.17
$ws.Range("I84").Value = 13709
$ws.Range("J84").Value = 248086
$ws.Range("K84").Value = 41127
$ws.Range("L84").Value = 744258
$ws.Range("M84").Value = -35823
$ws.Range("N84").Value = -754866
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2705.611
$ws.Range("J31").Value = 4101.3335
$ws.Range("L31").Value = 4101.3335
$ws.Range("N31").Value = -4691.3335
$ws.Range("H34").Value = 2705.611
$ws.Range("J34").Value = 4101.3335
$ws.Range("L34").Value = 4101.3335
$ws.Range("N34").Value = -4505.3335
$ws.Range("H94").Value = 3472.625
$ws.Range("I94").Value = 3316
$ws.Range("J94").Value = 3566.6
$ws.Range("K94").Value = 3316
$ws.Range("L94").Value = 3566.6
$ws.Range("M94").Value = -2865
$ws.Range("N94").Value = -4468.6
$ws.Range("H99").Value = 3187.8333
$ws.Range("I99").Value = 1950
$ws.Range("J99").Value = 3435.4
$ws.Range("K99").Value = 1950
$ws.Range("L99").Value = 3435.4
$ws.Range("M99").Value = -452
$ws.Range("N99").Value = -6431.4
$ws.Range("H122").Value = 1634.6
$ws.Range("I122").Value = 1471.7778
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 4415.3334
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -1965.3334
$ws.Range("N122").Value = -14200
$ws.Range("H126").Value = 3187.8333
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 3435.4
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 10306.2
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -15246.2
$ws.Range("H132").Value = 4308.5
$ws.Range("I132").Value = 4308.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12925.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -10395.5
$ws.Range("H134").Value = 32045
$ws.Range("I134").Value = 37912.703
$ws.Range("J134").Value = 5640.3335
$ws.Range("K134").Value = 113738.109
$ws.Range("L134").Value = 16921.0005
$ws.Range("M134").Value = -111203.109
$ws.Range("N134").Value = -21991.0005
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws.Range("H141").Value = 629997.6
$ws.Range("J141").Value = 629997.6
$ws.Range("L141").Value = 629997.6
$ws.Range("N141").Value = -640357.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 103
$ws.Range("I2").Value = 46.07143
$ws.Range("J2").Value = 235.83333
$ws.Range("K2").Value = 276.42858
$ws.Range("L2").Value = 1414.99998
$ws.Range("M2").Value = -163.42858
$ws.Range("N2").Value = -1640.99998
$ws.Range("H8").Value = 1197.6666
$ws.Range("I8").Value = 1197.6666
$ws.Range("K8").Value = 3592.9998
$ws.Range("M8").Value = -3453.9998
$ws.Range("H96").Value = 20007.666
$ws.Range("J96").Value = 24999
$ws.Range("L96").Value = 74997
$ws.Range("N96").Value = -79115
$ws.Range("H131").Value = 2508155.5
$ws.Range("J131").Value = 3232595
$ws.Range("L131").Value = 9697785
$ws.Range("N131").Value = -9707865
$ws.Range("H140").Value = 2763.1
$ws.Range("I140").Value = 2091.5625
$ws.Range("K140").Value = 6274.6875
$ws.Range("M140").Value = -1094.6875
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 971.5238000000001
$ws.Range("I97").Value = 713.9375
$ws.Range("J97").Value = 1795.8
$ws.Range("K97").Value = 713.9375
$ws.Range("L97").Value = 1795.8
$ws.Range("M97").Value = -217.9375
$ws.Range("N97").Value = -2787.8
$ws.Range("H122").Value = 1976.9412
$ws.Range("I122").Value = 1788.0625
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 5364.1875
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -2914.1875
$ws.Range("N122").Value = -19897
$ws.Range("H126").Value = 7471.5713
$ws.Range("I126").Value = 6319.5
$ws.Range("K126").Value = 18958.5
$ws.Range("M126").Value = -16488.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4693.9414
$ws.Range("I7").Value = 4130.3
$ws.Range("K7").Value = 4130.3
$ws.Range("M7").Value = -4018.3
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0
$ws.Range("H108").Value = 45000
$ws.Range("J108").Value = 45000
$ws.Range("L108").Value = 45000
$ws.Range("N108").Value = -52680
$ws.Range("H126").Value = 4693.9414
$ws.Range("I126").Value = 4130.3
$ws.Range("K126").Value = 12390.9
$ws.Range("M126").Value = -9920.900000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1286
$ws.Range("I96").Value = 1727.8
$ws.Range("J96").Value = 549.6667
$ws.Range("K96").Value = 1727.8
$ws.Range("L96").Value = 549.6667
$ws.Range("M96").Value = -354.8
$ws.Range("N96").Value = -3295.6667
$ws.Range("H126").Value = 43393.965
$ws.Range("I126").Value = 53687.76
$ws.Range("K126").Value = 161063.28
$ws.Range("M126").Value = -158593.28
$ws.Range("H132").Value = 29586.814
$ws.Range("I132").Value = 32175.016
$ws.Range("J132").Value = 1979.3334
$ws.Range("K132").Value = 96525.048
$ws.Range("L132").Value = 5938.0002
$ws.Range("M132").Value = -93995.048
$ws.Range("N132").Value = -10998.0002
